$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 2632
$ws.Range("E2").Value = 173
$ws.Range("F2").Value = 173
$ws.Range("G2").Value = 168
$ws.Range("H2").Value = 143
$ws.Range("I2").Value = 143
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3819
$ws.Range("L2").Value = 1628
$ws.Range("M2").Value = 2191
$ws.Range("N2").Value = 2187
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 145
$ws.Range("Q2").Value = 344
$ws.Range("R2").Value = 105
$ws.Range("S2").Value = -281
$ws.Range("T2").Value = 131
$ws.Range("U2").Value = 214
$ws.Range("V2").Value = 408
$ws.Range("W2").Value = 6.59
$ws.Range("X2").Value = 5.44
$ws.Range("Y2").Value = 6.7
$ws.Range("Z2").Value = 3.85
$ws.Range("AA2").Value = 74.3
$ws.Range("AB2").Value = 1379.89
$ws.Range("AC2").Value = 492
$ws.Range("AD2").Value = 10.4
$ws.Range("AE2").Value = 8094
$ws.Range("AF2").Value = 0.63
$ws.Range("AG2").Value = 80
$ws.Range("AH2").Value = 1.56
$ws.Range("AI2").Value = 15.09
$ws.Range("AJ2").Value = 29084925

# Row 3 updates
$ws.Range("D3").Value = 2695
$ws.Range("E3").Value = 212
$ws.Range("F3").Value = 212
$ws.Range("G3").Value = 220
$ws.Range("H3").Value = 164
$ws.Range("I3").Value = 164
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4228
$ws.Range("L3").Value = 1908
$ws.Range("M3").Value = 2320
$ws.Range("N3").Value = 2320
$ws.Range("P3").Value = 145
$ws.Range("Q3").Value = 103
$ws.Range("R3").Value = -160
$ws.Range("S3").Value = 40
$ws.Range("T3").Value = 96
$ws.Range("U3").Value = 7
$ws.Range("V3").Value = 532
$ws.Range("W3").Value = 7.88
$ws.Range("X3").Value = 6.1
$ws.Range("Y3").Value = 7.29
$ws.Range("Z3").Value = 4.09
$ws.Range("AA3").Value = 82.25
$ws.Range("AB3").Value = 1481.03
$ws.Range("AC3").Value = 565
$ws.Range("AD3").Value = 9.91
$ws.Range("AE3").Value = 8571
$ws.Range("AF3").Value = 0.65
$ws.Range("AG3").Value = 120
$ws.Range("AH3").Value = 2.14
$ws.Range("AI3").Value = 19.76
$ws.Range("AJ3").Value = 29084925
$ws.Range("O3").ClearContents()

# Row 4 updates
$ws.Range("D4").Value = 2671
$ws.Range("E4").Value = 236
$ws.Range("F4").Value = 236
$ws.Range("G4").Value = 217
$ws.Range("H4").Value = 162
$ws.Range("I4").Value = 162
$ws.Range("K4").Value = 4179
$ws.Range("L4").Value = 1774
$ws.Range("M4").Value = 2404
$ws.Range("N4").Value = 2404
$ws.Range("P4").Value = 145
$ws.Range("Q4").Value = 352
$ws.Range("R4").Value = -45
$ws.Range("S4").Value = -45
$ws.Range("T4").Value = 37
$ws.Range("U4").Value = 316
$ws.Range("V4").Value = 548
$ws.Range("W4").Value = 8.82
$ws.Range("X4").Value = 6.06
$ws.Range("Y4").Value = 6.85
$ws.Range("Z4").Value = 3.85
$ws.Range("AA4").Value = 73.81
$ws.Range("AB4").Value = 1569.6
$ws.Range("AC4").Value = 556
$ws.Range("AD4").Value = 8.52
$ws.Range("AE4").Value = 9075
$ws.Range("AF4").Value = 0.52
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 2.11
$ws.Range("AI4").Value = 16.38
$ws.Range("AJ4").Value = 29084925
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5 updates
$ws.Range("D5").Value = 2306
$ws.Range("E5").Value = 117
$ws.Range("F5").Value = 117
$ws.Range("G5").Value = 96
$ws.Range("H5").Value = 129
$ws.Range("I5").Value = 129
$ws.Range("K5").Value = 3812
$ws.Range("L5").Value = 1403
$ws.Range("M5").Value = 2408
$ws.Range("N5").Value = 2408
$ws.Range("P5").Value = 145
$ws.Range("Q5").Value = -95
$ws.Range("R5").Value = -4
$ws.Range("S5").Value = -103
$ws.Range("T5").Value = 14
$ws.Range("U5").Value = -110
$ws.Range("V5").Value = 552
$ws.Range("W5").Value = 5.07
$ws.Range("X5").Value = 5.59
$ws.Range("Y5").Value = 5.36
$ws.Range("Z5").Value = 3.23
$ws.Range("AA5").Value = 58.27
$ws.Range("AB5").Value = 1640.32
$ws.Range("AC5").Value = 443
$ws.Range("AD5").Value = 9.31
$ws.Range("AE5").Value = 9684
$ws.Range("AF5").Value = 0.43
$ws.Range("AG5").Value = 80
$ws.Range("AH5").Value = 1.94
$ws.Range("AI5").Value = 15.43
$ws.Range("AJ5").Value = 29084925
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6 updates
$ws.Range("D6").Value = 3062
$ws.Range("E6").Value = 102
$ws.Range("F6").Value = 102
$ws.Range("G6").Value = 715
$ws.Range("H6").Value = 692
$ws.Range("I6").Value = 692
$ws.Range("K6").Value = 6818
$ws.Range("L6").Value = 2467
$ws.Range("M6").Value = 4350
$ws.Range("N6").Value = 4350
$ws.Range("P6").Value = 327
$ws.Range("Q6").Value = 374
$ws.Range("R6").Value = 344
$ws.Range("S6").Value = 100
$ws.Range("T6").Value = 13
$ws.Range("U6").Value = 361
$ws.Range("V6").Value = 1134
$ws.Range("W6").Value = 3.33
$ws.Range("X6").Value = 22.6
$ws.Range("Y6").Value = 20.48
$ws.Range("Z6").Value = 13.02
$ws.Range("AA6").Value = 56.71
$ws.Range("AB6").Value = 1280.82
$ws.Range("AC6").Value = 1463
$ws.Range("AD6").Value = 1.67
$ws.Range("AE6").Value = 7185
$ws.Range("AF6").Value = 0.34
$ws.Range("AG6").Value = 60
$ws.Range("AH6").Value = 2.46
$ws.Range("AI6").Value = 5.25
$ws.Range("AJ6").Value = 65429516

# Row 7 updates
$ws.Range("D7:AI7").ClearContents()

# Row 8 updates
$ws.Range("D8:AI8").ClearContents()

# Row 9 updates
$ws.Range("D9:AI9").ClearContents()

